$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "59.562.12"
$ws.Range("E2").Value = "  +8.59%  "

Set-TextValue $ws.Range("D3") "2.572.56"
$ws.Range("E3").Value = "  +10.71%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "504.00"
$ws.Range("E5").Value = "  +6.77%  "

Set-TextValue $ws.Range("D6") "155.69"
$ws.Range("E6").Value = "  +8.83%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D8") "0.615"
$ws.Range("E8").Value = "  +4.64%  "

Set-TextValue $ws.Range("D9") "2.571.71"
$ws.Range("E9").Value = "  +11.14%  "

Set-TextValue $ws.Range("D10") "6.12"
$ws.Range("E10").Value = "  +13.48%  "

$ws.Range("E11").Value = "  +7.87%  "

$ws.Range("E12").Value = "  +6.88%  "

$ws.Range("E13").Value = "  +2.01%  "

Set-TextValue $ws.Range("D14") "3.020.71"
$ws.Range("E14").Value = "  +10.06%  "

Set-TextValue $ws.Range("D15") "59.402.41"
$ws.Range("E15").Value = "  +8.24%  "

Set-TextValue $ws.Range("D16") "21.78"
$ws.Range("E16").Value = "  +10.05%  "

$ws.Range("E17").Value = "  +6.86%  "

Set-TextValue $ws.Range("D18") "2.576.30"
$ws.Range("E18").Value = "  +10.17%  "

Set-TextValue $ws.Range("D19") "4.75"
$ws.Range("E19").Value = "  +5.10%  "

Set-TextValue $ws.Range("D20") "336.32"
$ws.Range("E20").Value = "  +7.72%  "

Set-TextValue $ws.Range("D21") "10.33"
$ws.Range("E21").Value = "  +8.85%  "

Set-TextValue $ws.Range("D22") "6.03"
$ws.Range("E22").Value = "  +8.78%  "

$ws.Range("E23").Value = "  +0.28%  "

Set-TextValue $ws.Range("D24") "59.94"
$ws.Range("E24").Value = "  +6.32%  "

$ws.Range("E25").Value = "  +6.88%  "

$ws.Range("E26").Value = "  +8.91%  "

Set-TextValue $ws.Range("D27") "2.684.47"
$ws.Range("E27").Value = "  +9.73%  "

Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  +0.09%  "

Set-TextValue $ws.Range("D29") "0.0₃0835"
$ws.Range("E29").Value = "  +13.50%  "

Set-TextValue $ws.Range("D30") "7.35"
$ws.Range("E30").Value = "  +4.66%  "

$ws.Range("E31").Value = "  -0.07%  "

Set-TextValue $ws.Range("D32") "157.28"
$ws.Range("E32").Value = "  +7.89%  "

Set-TextValue $ws.Range("D33") "19.37"
$ws.Range("E33").Value = "  +6.73%  "

$ws.Range("E34").Value = "  +7.26%  "

Set-TextValue $ws.Range("D35") "5.51"
$ws.Range("E35").Value = "  +9.83%  "

$ws.Range("E36").Value = "  +11.43%  "

$ws.Range("E37").Value = "  +10.44%  "

$ws.Range("E38").Value = "  +6.34%  "

Set-TextValue $ws.Range("D39") "3.74"
$ws.Range("E39").Value = "  +11.58%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D40") "1.44"
$ws.Range("E40").Value = "  +9.45%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D41") "295.09"
$ws.Range("E41").Value = "  +18.86%  "

$ws.Range("E42").Value = "  +4.95%  "

Set-TextValue $ws.Range("D43") "0.0573"
$ws.Range("E43").Value = "  +11.64%  "

$ws.Range("E44").Value = "  +4.10%  "

Set-TextValue $ws.Range("D45") "0.627"
$ws.Range("E45").Value = "  +10.49%  "

$ws.Range("E46").Value = "  +27.47%  "

$ws.Range("E47").Value = "  +0.26%  "

Set-TextValue $ws.Range("D48") "4.89"
$ws.Range("E48").Value = "  +13.50%  "

Set-TextValue $ws.Range("D49") "19.08"
$ws.Range("E49").Value = "  +15.66%  "

$ws.Range("E50").Value = "  +7.67%  "

Set-TextValue $ws.Range("D51") "10.26"
$ws.Range("E51").Value = "  +0.92%  "
